# Added ifoCAST full series evaluation.
# The error-table rows for quarters Q0..Q9 shift up by one evaluation window:
# each row now holds the figures that used to belong to the following row,
# the N column (G) decrements by one accordingly, and the final row (Q9,
# now N=5) receives newly computed statistics for the extended series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1872890445127373
$ws.Range("C2").Value = 0.5754152234511226
$ws.Range("D2").Value = 0.6110775347743189
$ws.Range("E2").Value = 0.7817144841784108
$ws.Range("F2").Value = 0.7875963862927191
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.1843166320775897
$ws.Range("C3").Value = 0.6273081474520239
$ws.Range("D3").Value = 0.5779702429309931
$ws.Range("E3").Value = 0.7602435418541831
$ws.Range("F3").Value = 0.7676788112505872
$ws.Range("G3").Value = 13

$ws.Range("B4").Value = 0.3221619536241352
$ws.Range("C4").Value = 0.5655998917628668
$ws.Range("D4").Value = 0.5382575468002283
$ws.Range("E4").Value = 0.7336603756509059
$ws.Range("F4").Value = 0.6884521947652315
$ws.Range("G4").Value = 12

$ws.Range("B5").Value = 0.2262967854219969
$ws.Range("C5").Value = 0.4712852108734478
$ws.Range("D5").Value = 0.2956388450555676
$ws.Range("E5").Value = 0.5437268110508876
$ws.Range("F5").Value = 0.5185281775945857
$ws.Range("G5").Value = 11

$ws.Range("B6").Value = 0.3318332233303713
$ws.Range("C6").Value = 0.399431789294712
$ws.Range("D6").Value = 0.2898093932768167
$ws.Range("E6").Value = 0.53833947772462
$ws.Range("F6").Value = 0.4468359196381605
$ws.Range("G6").Value = 10

$ws.Range("B7").Value = 0.3032448163457492
$ws.Range("C7").Value = 0.4442789123723634
$ws.Range("D7").Value = 0.3348453944406773
$ws.Range("E7").Value = 0.5786582708651776
$ws.Range("F7").Value = 0.5227322189947009
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = 0.3395226469355723
$ws.Range("C8").Value = 0.4954646948149822
$ws.Range("D8").Value = 0.3955394669362153
$ws.Range("E8").Value = 0.6289192849135852
$ws.Range("F8").Value = 0.5659518799878692
$ws.Range("G8").Value = 8

$ws.Range("B9").Value = 0.2802560178301312
$ws.Range("C9").Value = 0.3574868270362884
$ws.Range("D9").Value = 0.1586788892857869
$ws.Range("E9").Value = 0.3983451886062977
$ws.Range("F9").Value = 0.3057635732311509
$ws.Range("G9").Value = 7

$ws.Range("B10").Value = 0.3351585394422735
$ws.Range("C10").Value = 0.4364782826832794
$ws.Range("D10").Value = 0.3756151435636499
$ws.Range("E10").Value = 0.612874492505317
$ws.Range("F10").Value = 0.5620860044540215
$ws.Range("G10").Value = 6

$ws.Range("B11").Value = 0.3270154541542939
$ws.Range("C11").Value = 0.5200289994386147
$ws.Range("D11").Value = 0.3676720554669373
$ws.Range("E11").Value = 0.6063596750006858
$ws.Range("F11").Value = 0.5708906946728048
$ws.Range("G11").Value = 5

Write-Output "Applied ifoCAST full series evaluation update to B2:G11"
